$wb = $excel.ActiveWorkbook

# This report is generated after a "handback": the localized files have come
# back in sync with en-US. For each of the two language sheets (zh-cn, de-de)
# we:
#   - flip the Status column (B) from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two real rows (2 and 3)
#   - populate "Latest Target File" (E) / "Latest Handback File" (F) with the
#     same file reference as the handoff columns (A / C), including the
#     matching hyperlink
#   - stamp "Latest Handback DateTime" (G) with the real handback timestamp

function Set-HandbackRow($ws, $row, $handbackDateTime) {
    $statusText = "Handed back: in sync with en-US"

    # Status
    $ws.Cells.Item($row, 2).Value = $statusText

    # Find the existing hyperlinks anchored on column A (source file) and
    # column C (handoff file) for this row so we can mirror them.
    $srcHyperlink = $null
    $handoffHyperlink = $null
    foreach ($hl in $ws.Hyperlinks) {
        $r = $hl.Range
        if ($r.Row -eq $row -and $r.Column -eq 1) {
            $srcHyperlink = $hl
        }
        if ($r.Row -eq $row -and $r.Column -eq 3) {
            $handoffHyperlink = $hl
        }
    }

    # Source file (A) -> mirrored into "Latest Target File" (E)
    $srcValue = $ws.Cells.Item($row, 1).Value2
    $srcAddress = $srcHyperlink.Address
    $srcDisplay = $srcHyperlink.TextToDisplay

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $srcValue
    $ws.Hyperlinks.Add($eCell, $srcAddress, "", "", $srcDisplay) | Out-Null
    $eCell.Font.Underline = 2
    $eCell.Font.Color = 15570276

    # Handoff file (C) -> mirrored into "Latest Handback File" (F)
    $handoffValue = $ws.Cells.Item($row, 3).Value2
    $handoffAddress = $handoffHyperlink.Address
    $handoffDisplay = $handoffHyperlink.TextToDisplay

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $handoffValue
    $ws.Hyperlinks.Add($fCell, $handoffAddress, "", "", $handoffDisplay) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = 15570276

    # Latest Handback DateTime (G)
    $ws.Cells.Item($row, 7).Value = $handbackDateTime
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZh 2 "2016-03-10 11:56:47"
Set-HandbackRow $wsZh 3 "2016-03-10 11:56:47"

$wsDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDe 2 "2016-03-10 11:56:53"
Set-HandbackRow $wsDe 3 "2016-03-10 11:56:53"
